$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.66950562919271
$ws.Range("C2").Value = -0.2895456687149427
$ws.Range("D2").Value = -0.6733414736251095
$ws.Range("E2").Value = 0.0115444171491989
$ws.Range("F2").Value = -0.6804809672324722
$ws.Range("G2").Value = -0.4410326232298434
$ws.Range("H2").Value = -0.4043046464706727
$ws.Range("I2").Value = 0.4221894279166428
$ws.Range("J2").Value = -0.2336087822049224
$ws.Range("K2").Value = -0.8917590219305432

$ws.Range("B3").Value = -0.3724265895266604
$ws.Range("C3").Value = 0.3124593012476481
$ws.Range("D3").Value = -0.3795660831340231
$ws.Range("E3").Value = -0.1401177391313942
$ws.Range("F3").Value = -0.1033897623722235
$ws.Range("G3").Value = 0.723104312015092
$ws.Range("H3").Value = 0.06730610189352677
$ws.Range("I3").Value = -0.5908441378320941
$ws.Range("J3").Value = -0.191729189980311
$ws.Range("K3").Value = -0.2531657321018884

$ws.Range("B4").Value = 0.2606382006276979
$ws.Range("C4").Value = 0.5000865446303268
$ws.Range("D4").Value = 0.5368145213894975
$ws.Range("E4").Value = 1.363308595776813
$ws.Range("F4").Value = 0.7075103856552477
$ws.Range("G4").Value = 0.0493601459296269
$ws.Range("H4").Value = 0.44847509378141
$ws.Range("I4").Value = 0.3870385516598326
$ws.Range("J4").Value = -0.04712738345727097
$ws.Range("K4").Value = 0.9316740086218059

$ws.Range("B5").Value = 0.4642269187270481
$ws.Range("C5").Value = 1.290720993114364
$ws.Range("D5").Value = 0.6349227829927984
$ws.Range("E5").Value = -0.02322745673282245
$ws.Range("F5").Value = 0.3758874911189606
$ws.Range("G5").Value = 0.3144509489973832
$ws.Range("H5").Value = -0.1197149861197203
$ws.Range("I5").Value = 0.8590864059593566
$ws.Range("J5").Value = 0.6319705015114304
$ws.Range("K5").Value = -0.001489227914658109

$ws.Range("B6").Value = 0.5792317386292632
$ws.Range("C6").Value = -0.07891850109635756
$ws.Range("D6").Value = 0.3201964467554255
$ws.Range("E6").Value = 0.2587599046338481
$ws.Range("F6").Value = -0.1754060304832554
$ws.Range("G6").Value = 0.8033953615958215
$ws.Range("H6").Value = 0.5762794571478953
$ws.Range("I6").Value = -0.05718027227819322
$ws.Range("J6").Value = 0.3702166863774111
$ws.Range("K6").Value = 0.2038797316892642

$ws.Range("B7").Value = -0.03386271341853997
$ws.Range("C7").Value = -0.09529925554011737
$ws.Range("D7").Value = -0.529465190657221
$ws.Range("E7").Value = 0.449336201421856
$ws.Range("F7").Value = 0.2222202969739298
$ws.Range("G7").Value = -0.4112394324521587
$ws.Range("H7").Value = 0.01615752620344563
$ws.Range("I7").Value = -0.1501794284847013
$ws.Range("J7").Value = -0.02746841204387546
$ws.Range("K7").Value = 0.01976532451126484

$ws.Range("B8").Value = -0.6197834074567948
$ws.Range("C8").Value = 0.3590179846222821
$ws.Range("D8").Value = 0.1319020801743559
$ws.Range("E8").Value = -0.5015576492517326
$ws.Range("F8").Value = -0.07416069059612829
$ws.Range("G8").Value = -0.2404976452842752
$ws.Range("H8").Value = -0.1177866288434494
$ws.Range("I8").Value = -0.07055289228830908
$ws.Range("J8").Value = -0.4671716238107607
$ws.Range("K8").Value = -0.2131195886870406

$ws.Range("B9").Value = 0.2002248142788162
$ws.Range("C9").Value = -0.4332349151472724
$ws.Range("D9").Value = -0.005837956491668017
$ws.Range("E9").Value = -0.1721749111798149
$ws.Range("F9").Value = -0.0494638947389891
$ws.Range("G9").Value = -0.002230158183848807
$ws.Range("H9").Value = -0.3988488897063004
$ws.Range("I9").Value = -0.1447968545825803
$ws.Range("J9").Value = 0.07532753529099229
$ws.Range("K9").Value = -0.1604609132956759

$ws.Range("B10").Value = -0.1044995345215733
$ws.Range("C10").Value = -0.2708364892097202
$ws.Range("D10").Value = -0.1481254727688944
$ws.Range("E10").Value = -0.1008917362137541
$ws.Range("F10").Value = -0.4975104677362057
$ws.Range("G10").Value = -0.2434584326124856
$ws.Range("H10").Value = -0.02333404273891299
$ws.Range("I10").Value = -0.2591224913255812
$ws.Range("J10").Value = -0.4380952487963659
$ws.Range("K10").Value = -0.4486557790556869

$ws.Range("B11").Value = -0.07111589452722158
$ws.Range("C11").Value = -0.02388215797208129
$ws.Range("D11").Value = -0.4205008894945329
$ws.Range("E11").Value = -0.1664488543708128
$ws.Range("F11").Value = 0.0536755355027598
$ws.Range("G11").Value = -0.1821129130839084
$ws.Range("H11").Value = -0.3610856705546931
$ws.Range("I11").Value = -0.3716462008140141
$ws.Range("J11").Value = -0.5793653109721442
$ws.Range("K11").Value = -0.3621564114972394

$ws.Range("B12").Value = -0.2489848374783721
$ws.Range("C12").Value = 0.005067197645347965
$ws.Range("D12").Value = 0.2251915875189206
$ws.Range("E12").Value = -0.0105968610677476
$ws.Range("F12").Value = -0.1895696185385323
$ws.Range("G12").Value = -0.2001301487978533
$ws.Range("H12").Value = -0.4078492589559834
$ws.Range("I12").Value = -0.1906403594810787
$ws.Range("J12").Value = 0.1190458097769828
$ws.Range("K12").Value = -0.03525832405502122

$ws.Range("B13").Value = 0.1877264853751708
$ws.Range("C13").Value = -0.04806196321149736
$ws.Range("D13").Value = -0.2270347206822821
$ws.Range("E13").Value = -0.2375952509416031
$ws.Range("F13").Value = -0.4453143610997332
$ws.Range("G13").Value = -0.2281054616248284
$ws.Range("H13").Value = 0.08158070763323305
$ws.Range("I13").Value = -0.07272342619877098
$ws.Range("J13").Value = 0.5777029950204122
$ws.Range("K13").Value = 0.361537946402203

$ws.Range("B14").Value = -0.2489372330903247
$ws.Range("C14").Value = -0.2594977633496457
$ws.Range("D14").Value = -0.4672168735077758
$ws.Range("E14").Value = -0.2500079740328711
$ws.Range("F14").Value = 0.05967819522519041
$ws.Range("G14").Value = -0.09462593860681362
$ws.Range("H14").Value = 0.5558004826123696
$ws.Range("I14").Value = 0.3396354339941604
$ws.Range("J14").Value = -0.0673936950407959
$ws.Range("K14").Value = 0.4316485221366239

$ws.Range("B15").Value = -0.3791109816618064
$ws.Range("C15").Value = -0.1619020821869017
$ws.Range("D15").Value = 0.1477840870711598
$ws.Range("E15").Value = -0.006520046760844223
$ws.Range("F15").Value = 0.643906374458339
$ws.Range("G15").Value = 0.4277413258401298
$ws.Range("H15").Value = 0.0207121968051735
$ws.Range("I15").Value = 0.5197544139825933
$ws.Range("J15").Value = 0.3609055008270807
$ws.Range("K15").Value = 0.4170693217973248

$ws.Range("B16").Value = 0.241497203637704
$ws.Range("C16").Value = 0.08719306980569996
$ws.Range("D16").Value = 0.7376194910248832
$ws.Range("E16").Value = 0.521454442406674
$ws.Range("F16").Value = 0.1144253133717177
$ws.Range("G16").Value = 0.6134675305491375
$ws.Range("H16").Value = 0.4546186173936249
$ws.Range("I16").Value = 0.5107824383638689
$ws.Range("J16").Value = 2.760585277975261
$ws.Range("K16").Value = 10.22973560883255

$ws.Range("B17").Value = 0.09539661607613537
$ws.Range("C17").Value = 0.7458230372953185
$ws.Range("D17").Value = 0.5296579886771094
$ws.Range("E17").Value = 0.1226288596421531
$ws.Range("F17").Value = 0.6216710768195729
$ws.Range("G17").Value = 0.4628221636640603
$ws.Range("H17").Value = 0.5189859846343043
$ws.Range("I17").Value = 2.768788824245696
$ws.Range("J17").Value = 10.23793915510299
$ws.Range("K17").Value = -7.933322645034479

$ws.Range("B18").Value = 0.7432334766438504
$ws.Range("C18").Value = 0.5270684280256412
$ws.Range("D18").Value = 0.1200392989906849
$ws.Range("E18").Value = 0.6190815161681047
$ws.Range("F18").Value = 0.4602326030125921
$ws.Range("G18").Value = 0.5163964239828361
$ws.Range("H18").Value = 2.766199263594229
$ws.Range("I18").Value = 10.23534959445152
$ws.Range("J18").Value = -7.935912205685947
$ws.Range("K18").Value = 0.2037593412932021

$ws.Range("B19").Value = 0.5336017606086931
$ws.Range("C19").Value = 0.1265726315737368
$ws.Range("D19").Value = 0.6256148487511566
$ws.Range("E19").Value = 0.466765935595644
$ws.Range("F19").Value = 0.522929756565888
$ws.Range("G19").Value = 2.77273259617728
$ws.Range("H19").Value = 10.24188292703457
$ws.Range("I19").Value = -7.929378873102896
$ws.Range("J19").Value = 0.2102926738762539
$ws.Range("K19").Value = 2.357863145675646

$ws.Range("B20").Value = 0.01535593659466611
$ws.Range("C20").Value = 0.5143981537720859
$ws.Range("D20").Value = 0.3555492406165733
$ws.Range("E20").Value = 0.4117130615868174
$ws.Range("F20").Value = 2.661515901198209
$ws.Range("G20").Value = 10.1306662320555
$ws.Range("H20").Value = -8.040595568081965
$ws.Range("I20").Value = 0.09907597889718328
$ws.Range("J20").Value = 2.246646450696576
$ws.Range("K20").Value = -1.193383535826754

$ws.Range("B21").Value = 0.4726520097851818
$ws.Range("C21").Value = 0.3138030966296693
$ws.Range("D21").Value = 0.3699669175999133
$ws.Range("E21").Value = 2.619769757211305
$ws.Range("F21").Value = 10.0889200880686
$ws.Range("G21").Value = -8.082341712068869
$ws.Range("H21").Value = 0.05732983491027921
$ws.Range("I21").Value = 2.204900306709672
$ws.Range("J21").Value = -1.235129679813658
$ws.Range("K21").Value = -1.296303858227859

$ws.Range("B22").Value = 0.3391181593623842
$ws.Range("C22").Value = 0.3952819803326282
$ws.Range("D22").Value = 2.64508481994402
$ws.Range("E22").Value = 10.11423515080131
$ws.Range("F22").Value = -8.057026649336155
$ws.Range("G22").Value = 0.0826448976429941
$ws.Range("H22").Value = 2.230215369442386
$ws.Range("I22").Value = -1.209814617080943
$ws.Range("J22").Value = -1.270988795495144
$ws.Range("K22").Value = 0.77046714166111

$ws.Range("B23").Value = 0.3993206453755421
$ws.Range("C23").Value = 2.649123484986935
$ws.Range("D23").Value = 10.11827381584423
$ws.Range("E23").Value = -8.052987984293241
$ws.Range("F23").Value = 0.08668356268590799
$ws.Range("G23").Value = 2.2342540344853
$ws.Range("H23").Value = -1.20577595203803
$ws.Range("I23").Value = -1.266950130452231
$ws.Range("J23").Value = 0.7745058067040239
$ws.Range("K23").Value = 0.1919112696703829

$ws.Range("B24").Value = 2.530246348295912
$ws.Range("C24").Value = 9.999396679153206
$ws.Range("D24").Value = -8.171865120984263
$ws.Range("E24").Value = -0.03219357400511441
$ws.Range("F24").Value = 2.115376897794278
$ws.Range("G24").Value = -1.324653088729052
$ws.Range("H24").Value = -1.385827267143253
$ws.Range("I24").Value = 0.6556286700130015
$ws.Range("J24").Value = 0.07303413297936051
$ws.Range("K24").Value = 0.1687091737887183

Write-Output "done"